$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = "Smooth Grad-CAM++ 내용 정리 [XAI-4]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/145"

$ws.Range("D23").Value = "failed to create cublas handle: CUBLAS_STATUS_ALLOC_FAILED"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2783"

$ws.Range("D37").Value = "[Paper Review] Don't Stop Pretraining: Adapt Language Models to Domains and Tasks"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1456&mod=document&pageid=1"

$ws.Range("D41").Value = "클라우드 환경에서의 데이터베이스 가용성"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/%ed%81%b4%eb%9d%bc%ec%9a%b0%eb%93%9c-%ed%99%98%ea%b2%bd%ec%97%90%ec%84%9c%ec%9d%98-%eb%8d%b0%ec%9d%b4%ed%84%b0%eb%b2%a0%ec%9d%b4%ec%8a%a4-%ea%b0%80%ec%9a%a9%ec%84%b1/"

$ws.Range("D51").Value = "[python] pytesseract 사용중 print 안 되는 애러 해결법"
$ws.Range("E51").Value = "https://bskyvision.com/1169"
